$wb = $excel.ActiveWorkbook

# --- Rename the "Include from NMDP Ethnicity C" tab to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from NMDP Ethnicity C")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Date: 2024-05-20T10:22:59-05:00 -> 2024-11-11T17:53:38-06:00
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows for the new row.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
